$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")
$ws.Activate()

# Insert a new row above the current row 4 (the blank separator before the
# "facilities" choice list), shifting all subsequent rows down by one.
$ws.Rows.Item(4).Insert()

# Carry over the blank-cell formatting that already exists a couple of rows
# down (columns D:Z of the facilities rows) onto the new row's trailing
# columns, so the new row looks consistent with its neighbours.
$ws.Range("D6:Z6").Copy()
$ws.Range("D4:Z4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New "Ndebele" language_preference choice, following the same pattern as
# the existing "english" (row 2) and "shona" (row 3) choices.
$ws.Range("A4").Value = "language_preference"
$ws.Range("B4").Value = "ndebele"
$ws.Range("C4").Value = "Ndebele"
